$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "a"
$ws.Range("C3").Value = "b"
$ws.Range("D3").Value = "c"
$ws.Range("E3").Value = "d"

$ws.Range("E3").Select()
